# Apply updated crypto price/volume figures to the worksheet.
# Column D ("Price") cells that look like plain numbers are written with a
# leading apostrophe so Excel stores them as text (matching the original
# inline-string cells) instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.200.36'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '2.600.92'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''582.74'
$ws.Range('E5').Value = '  +3.20%  '
$ws.Range('D6').Value = '''142.98'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').Value = '''6.55'
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').Value = '''0.372'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').Value = '3.062.32'
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('D14').Value = '''24.41'
$ws.Range('D15').Value = '60.195.91'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '2.604.48'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '''11.36'
$ws.Range('E18').Value = '  +3.88%  '
$ws.Range('D19').Value = '''4.62'
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('D20').Value = '''345.68'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').Value = '''6.90'
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('E23').Value = '  +2.62%  '
$ws.Range('D24').Value = '''63.67'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('D27').Value = '''7.99'
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('E28').Value = '  +8.95%  '
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('D32').Value = '''166.97'
$ws.Range('E32').Value = '  +4.83%  '
$ws.Range('D33').Value = '''19.43'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').Value = '''1.31'
$ws.Range('E34').Value = '  +9.72%  '
$ws.Range('D35').Value = '''4.25'
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('D36').Value = '''0.982'
$ws.Range('E36').Value = '  +2.56%  '
$ws.Range('E37').Value = '  +4.26%  '
$ws.Range('D38').Value = '''38.15'
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('D39').Value = '''313.16'
$ws.Range('E39').Value = '  +3.62%  '
$ws.Range('E40').Value = '  +1.47%  '
$ws.Range('D41').Value = '''0.841'
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('D42').Value = '''135.50'
$ws.Range('E42').Value = '  -3.88%  '
$ws.Range('D43').Value = '''0.0994'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').Value = '''0.0549'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('E49').Value = '  +3.03%  '
$ws.Range('D50').Value = '''19.88'
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('E51').Value = '  +0.43%  '
